# Word COM-interop script applying the resume content updates described by
# the commit "resume and doc updates".
#
# Strategy: for plain text swaps we replace the *entire* containing
# paragraph (including its paragraph mark) with freshly authored
# WordprocessingML via Range.InsertXML - this keeps xml:space="preserve"
# on the text runs (Find.Execute's Replace parameter silently drops that
# attribute) and leaves paragraph/run formatting untouched. For the new
# bullet points / Skills line we insert whole new <w:p> paragraphs next to
# an anchor paragraph, also via InsertXML.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParagraphByAnchor($anchorText) {
    $range = $d.Content
    $ok = $range.Find.Execute($anchorText)
    if (-not $ok) {
        throw "Find.Execute failed to find: $anchorText"
    }
    return $range.Paragraphs(1)
}

function Replace-FirstParagraphText($anchorText, $newText) {
    # Swaps a single-run "FirstParagraph"-styled paragraph's text, keeping
    # style + xml:space="preserve" intact.
    $para = Get-ParagraphByAnchor $anchorText
    $xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"FirstParagraph`"/></w:pPr><w:r><w:t xml:space=`"preserve`">$newText</w:t></w:r></w:p>"
    $para.Range.InsertXML($xml)
}

function Insert-ParagraphXmlBefore($anchorText, $xmlFragment) {
    # Inserts the given WordprocessingML paragraph(s) immediately before the
    # paragraph whose text contains $anchorText.
    $para = Get-ParagraphByAnchor $anchorText
    $startPos = $para.Range.Start
    $ins = $d.Range($startPos, $startPos)
    $ins.InsertXML($xmlFragment)
}

function Insert-ParagraphXmlAfter($anchorText, $xmlFragment) {
    # Inserts the given WordprocessingML paragraph(s) immediately after the
    # paragraph whose text contains $anchorText (targeting the position of
    # its paragraph mark so any bookmarkEnd/bookmarkStart that immediately
    # follows stays after the newly-inserted paragraph).
    $para = Get-ParagraphByAnchor $anchorText
    $insPos = $para.Range.End - 1
    $ins = $d.Range($insPos, $insPos)
    $ins.InsertXML($xmlFragment)
}

# ---------------------------------------------------------------------
# 1. Dfinity: rewrite the company summary paragraph.
# ---------------------------------------------------------------------
Replace-FirstParagraphText `
    "Blockchain company developing the Internet Computer, a decentralized platform for running scalable smart contracts and dApps without traditional cloud infrastructure." `
    "DFINITY Foundation is a not-for-profit organization developing the Internet Computer, a revolutionary blockchain network that transforms the internet into a decentralized cloud. Home to Caffeine AI, the world’s first self-writing apps platform that enables users to create fully decentralized applications using natural language conversations, deployed directly on the blockchain without traditional coding."

# ---------------------------------------------------------------------
# 2. Dfinity: insert four new bullet paragraphs before the existing
#    "Engineered 130 website frontend templates..." bullet.
# ---------------------------------------------------------------------
$dfinityBullets = @"
<w:p $wNs><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1001"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Shipped</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">production-ready applications end-to-end for Caffeine’s launch, including Tiny Tasks (featured in the live launch demo) and Habitual, demonstrating full ownership from concept to deployment.</w:t></w:r></w:p><w:p $wNs><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1001"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Led</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">platform-wide UI/UX improvements across the App Marketplace, standardized theming architecture, and resolved critical user experience issues like Internet Identity sign-in flows.</w:t></w:r></w:p><w:p $wNs><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1001"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Enabled</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">cross-team collaboration during Caffeine’s high-stakes product launch through comprehensive technical documentation, shared component libraries, and developer enablement tools.</w:t></w:r></w:p><w:p $wNs><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1001"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Architected</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">modern React/TypeScript solutions with mobile-first design, improving both developer experience and end-user experience while maintaining high delivery pace.</w:t></w:r></w:p>
"@
Insert-ParagraphXmlBefore "130 website frontend templates to be used as a RAG for our Caffeine AI Website Builder." $dfinityBullets

# ---------------------------------------------------------------------
# 3. Dfinity: insert a new "Skills:" paragraph right after the "Updated
#    the styling for the KYC-site." bullet (still inside the Dfinity
#    bookmark range).
# ---------------------------------------------------------------------
$dfinitySkills = @"
<w:p $wNs><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Skills:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">TypeScript · React · Python · CrewAI · MCP Servers · Tailwind · Technical Documentation · Mentorship</w:t></w:r></w:p>
"@
Insert-ParagraphXmlAfter "the styling for the KYC-site." $dfinitySkills

# ---------------------------------------------------------------------
# 4. Bill: rewrite the company summary paragraph.
# ---------------------------------------------------------------------
Replace-FirstParagraphText `
    "Payment platform for small and medium businesses. Worked with Angular & RxJS on the revenue-generating International Payments Team." `
    "Automated financial operations platform providing end-to-end bill payments, invoicing, and accounts payable/receivable management for SMBs and midsize companies. AI-enabled platform that streamlines cash flow and payment processes. Worked with Angular &amp; RxJS on the revenue-generating International Payments Team."

# ---------------------------------------------------------------------
# 5. Williams Sonoma Inc: rewrite the company summary paragraph.
# ---------------------------------------------------------------------
Replace-FirstParagraphText `
    "E-Commerce sites for housewares and home decor. Worked on the Micro Front End team, transitioning Production pages from 6 brands, into a unified Vue project." `
    "Premier multi-channel specialty retailer of high-quality home furnishings and kitchenware, operating 625+ stores globally and distributing to 60+ countries. Portfolio includes Williams Sonoma, Pottery Barn, Pottery Barn Kids, PBteen, West Elm, Williams-Sonoma Home, Mark and Graham, and Rejuvenation brands. Worked on the Micro Front End team, transitioning Production pages from 6 brands, into a unified Vue project."

# ---------------------------------------------------------------------
# 6. Williams Sonoma Inc: append " · Mentorship" to the Skills line
#    (3rd run of that paragraph holds the skill list text).
# ---------------------------------------------------------------------
$wsSkillsPara = Get-ParagraphByAnchor "Javascript · Vue · SASS · NodeJS · Shell · Markdown · Adobe Experience Manager"
$wsSkillsXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"FirstParagraph`"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space=`"preserve`">Skills:</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t xml:space=`"preserve`">Javascript · Vue · SASS · NodeJS · Shell · Markdown · Adobe Experience Manager · Mentorship</w:t></w:r></w:p>"
$wsSkillsPara.Range.InsertXML($wsSkillsXml)

# ---------------------------------------------------------------------
# 7. NewlyWords: rewrite the company summary paragraph.
# ---------------------------------------------------------------------
Replace-FirstParagraphText `
    "A memory book platform, where friends & family collaborate on a commemorative digital or physical book." `
    "Collaborative memory book platform that enables users to collect letters, photos, and memories from friends, family, and coworkers for special occasions like retirements, birthdays, and anniversaries. Contributors submit content through customizable templates, which are then compiled into professionally printed hardcover books with immediate PDF download options."

# ---------------------------------------------------------------------
# 8. Hawaii United Okinawa Association: rewrite the company summary
#    paragraph.
# ---------------------------------------------------------------------
Replace-FirstParagraphText `
    "A philanthropic organization dedicated to cultural exchange between Hawaii & Japan." `
    "Non-profit organization serving to promote, perpetuate, and preserve Okinawan culture in Hawaii. Represents 50 member clubs with combined membership exceeding 40,000, publishing the bi-monthly Uchinanchu Newsletter and commemorating 125 years of Okinawan immigration to Hawaii and emigration to the world."

# ---------------------------------------------------------------------
# 9. LifeStreet: rewrite the company summary paragraph.
# ---------------------------------------------------------------------
Replace-FirstParagraphText `
    "Mobile marketing platform. Engineered & Designed mobile marketing materials." `
    "Mobile-first demand-side platform (DSP) specializing in programmatic advertising and user acquisition for mobile apps and websites. Features Nero platform with true bidding transparency, predictive analytics, and ROAS optimization for performance advertisers seeking scalable campaign management."

Write-Output "All edits applied"
